# Add a missing entry for VehId = 439.
# The row belongs right before the existing VehId=440 row (row 239), so insert
# a new row there, pushing all subsequent rows down by one, then populate it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Rows.Item(239).Insert()

$ws.Range("A239").Value = 439
$ws.Range("B239").Value = "HEV"
$ws.Range("C239").Value = "NO DATA"
$ws.Range("D239").Value = "4-GAS/ELECTRIC 2.0L"
$ws.Range("E239").Value = "CVT"
$ws.Range("F239").Value = "NO DATA"
$ws.Range("G239").Value = 3500

[void]$ws.Range("C6").Select()
